$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "57.962.52"
Set-TextValue "E2" "  -3.05%  "
Set-TextValue "D3" "2.299.06"
Set-TextValue "E3" "  -3.37%  "
Set-TextValue "D4" "0.997"
Set-TextValue "E4" "  -0.32%  "
Set-TextValue "D5" "534.04"
Set-TextValue "E5" "  -3.91%  "
Set-TextValue "D6" "130.74"
Set-TextValue "E6" "  -2.00%  "
Set-TextValue "D7" "0.998"
Set-TextValue "E7" "  -0.17%  "
Set-TextValue "D8" "0.580"
Set-TextValue "E8" "  -1.19%  "
Set-TextValue "D9" "2.281.32"
Set-TextValue "E9" "  -4.09%  "
Set-TextValue "D10" "0.0996"
Set-TextValue "E10" "  -4.99%  "
Set-TextValue "D11" "5.42"
Set-TextValue "E11" "  -3.52%  "
Set-TextValue "E12" "  -0.36%  "
Set-TextValue "D13" "0.330"
Set-TextValue "E13" "  -4.04%  "
Set-TextValue "D14" "23.44"
Set-TextValue "E14" "  -3.93%  "
Set-TextValue "D15" "2.684.15"
Set-TextValue "E15" "  -4.27%  "
Set-TextValue "D16" "57.764.15"
Set-TextValue "E16" "  -3.27%  "
Set-TextValue "E17" "  -3.91%  "
Set-TextValue "D18" "2.261.96"
Set-TextValue "E18" "  -4.30%  "
Set-TextValue "D19" "10.51"
Set-TextValue "E19" "  -5.44%  "
Set-TextValue "D20" "4.21"
Set-TextValue "E20" "  -5.81%  "
Set-TextValue "D21" "313.47"
Set-TextValue "E21" "  -1.92%  "
Set-TextValue "D22" "6.36"
Set-TextValue "E22" "  -4.97%  "
Set-TextValue "E23" "  -0.03%  "
Set-TextValue "D24" "62.50"
Set-TextValue "E24" "  -2.56%  "
Set-TextValue "D25" "0.166"
Set-TextValue "E25" "  -4.07%  "
Set-TextValue "D26" "0.997"
Set-TextValue "E26" "  -0.60%  "
Set-TextValue "D27" "8.01"
Set-TextValue "E27" "  -4.53%  "
Set-TextValue "D28" "1.29"
Set-TextValue "E28" "  -5.61%  "
Set-TextValue "D29" "170.48"
Set-TextValue "E29" "  +0.51%  "
Set-TextValue "D30" "1.70"
Set-TextValue "E30" "  -5.51%  "
Set-TextValue "D31" "0.0₃0719"
Set-TextValue "E31" "  -5.06%  "
Set-TextValue "D32" "5.77"
Set-TextValue "E32" "  -4.54%  "
Set-TextValue "D33" "1.04"
Set-TextValue "E33" "  -4.21%  "
Set-TextValue "B34" "PolygonEcosystemToken"
Set-TextValue "C34" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D34" "0.376"
Set-TextValue "E34" "  -4.71%  "
Set-TextValue "B35" "USDe"
Set-TextValue "C35" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D35" "0.998"
Set-TextValue "E35" "  -0.04%  "
Set-TextValue "D36" "17.77"
Set-TextValue "E36" "  -2.23%  "
Set-TextValue "D37" "0.999"
Set-TextValue "E38" "  -6.53%  "
Set-TextValue "D39" "3.89"
Set-TextValue "E39" "  -5.86%  "
Set-TextValue "D40" "37.98"
Set-TextValue "E40" "  -1.62%  "
Set-TextValue "D41" "1.49"
Set-TextValue "E41" "  -5.35%  "
Set-TextValue "D42" "140.90"
Set-TextValue "E42" "  -3.44%  "
Set-TextValue "D43" "287.42"
Set-TextValue "E43" "  -9.73%  "
Set-TextValue "D44" "3.41"
Set-TextValue "E44" "  -3.26%  "
Set-TextValue "D45" "0.0947"
Set-TextValue "E45" "  -2.33%  "
Set-TextValue "D46" "0.0496"
Set-TextValue "E46" "  -2.86%  "
Set-TextValue "D47" "0.555"
Set-TextValue "E47" "  -2.81%  "
Set-TextValue "D48" "18.04"
Set-TextValue "E48" "  -8.25%  "
Set-TextValue "D49" "0.0210"
Set-TextValue "E49" "  -3.21%  "
Set-TextValue "D50" "10.92"
Set-TextValue "E50" "  -1.14%  "
Set-TextValue "E51" "  +82.89%  "
